$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '89.810.78'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.037.27'
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.56'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '612.72'
$ws.Range('E6').Value = '  -3.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.366'
$ws.Range('E7').Value = '  -7.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.896'
$ws.Range('E8').Value = '  +14.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.032.08'
$ws.Range('E10').Value = '  -3.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.667'
$ws.Range('E11').Value = '  +18.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.189'
$ws.Range('E12').Value = '  +5.54%  '
$ws.Range('E13').Value = '  -4.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.32'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.372.66'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.40'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.586.77'
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.022.68'
$ws.Range('E18').Value = '  -3.51%  '
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000219'
$ws.Range('E20').Value = '  -3.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.43'
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '424.66'
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.27'
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.04'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '82.93'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.56'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.23'
$ws.Range('E29').Value = '  +22.95%  '
$ws.Range('E30').Value = '  +1.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.52'
$ws.Range('E31').Value = '  +4.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.74'
$ws.Range('E32').Value = '  -6.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '501.67'
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.65'
$ws.Range('E34').Value = '  -4.62%  '
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.77'
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  -3.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.133'
$ws.Range('E38').Value = '  -9.15%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +3.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.359'
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '143.15'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.58'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0694'
$ws.Range('E47').Value = '  +6.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.23'
$ws.Range('E48').Value = '  +7.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '160.90'
$ws.Range('E49').Value = '  -2.28%  '
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.592'
$ws.Range('E51').Value = '  -1.38%  '
